# Update "想去人数" (want-to-go count) values for two events that appear
# on both the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): row 3 -> 熊喵M动漫嘉年华·万圣派对, row 4 -> 万圣漫控嘉年华10
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 153
$wsExpo.Range("F4").Value = 716

# Sheet "全部类型" (all types): row 4 -> 熊喵M动漫嘉年华·万圣派对, row 5 -> 万圣漫控嘉年华10
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 153
$wsAll.Range("F5").Value = 716
